# Updated cryptos list (prices + 1h volume deltas) on the "cryptos" sheet.
# Generated from the authoritative before/after cell diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds plain-text numeric-looking price strings (e.g. "303.20").
# Excel auto-coerces such text to a float on assignment, which would drop
# trailing zeros / precision. Mark those specific cells as Text first so the
# literal string survives, matching the source data (all cells are text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '45.721.04'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.578.16'
$ws.Range('E3').Value = '  +8.30%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '303.20'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '98.41'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').Value = '0.587'
$ws.Range('E7').Value = '  +3.91%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.563'
$ws.Range('E9').Value = '  +10.40%  '
$ws.Range('D10').Value = '37.74'
$ws.Range('E10').Value = '  +9.09%  '
$ws.Range('D11').Value = '0.0823'
$ws.Range('E11').Value = '  +4.35%  '
$ws.Range('D12').Value = '7.91'
$ws.Range('E12').Value = '  +10.77%  '
$ws.Range('D13').Value = '2.987.17'
$ws.Range('E13').Value = '  +8.82%  '
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '2.617.48'
$ws.Range('E15').Value = '  +9.47%  '
$ws.Range('D16').Value = '0.885'
$ws.Range('E16').Value = '  +6.84%  '
$ws.Range('D17').Value = '14.64'
$ws.Range('E17').Value = '  +6.55%  '
$ws.Range('D18').Value = '45.930.22'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '12.90'
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').Value = '0.0₃0992'
$ws.Range('E20').Value = '  +4.23%  '
$ws.Range('D21').Value = '6.57'
$ws.Range('E21').Value = '  +8.57%  '
$ws.Range('D22').Value = '70.08'
$ws.Range('E22').Value = '  +4.70%  '
$ws.Range('D23').Value = '251.07'
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  +4.94%  '
$ws.Range('D25').Value = '2.19'
$ws.Range('E25').Value = '  +13.51%  '
$ws.Range('D26').Value = '27.15'
$ws.Range('E26').Value = '  +29.67%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '10.31'
$ws.Range('E28').Value = '  +5.58%  '
$ws.Range('D29').Value = '39.01'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').Value = '2.24'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('B31').Value = 'LidoDAOToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').Value = '3.66'
$ws.Range('E31').Value = '  -2.95%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.99'
$ws.Range('E32').Value = '  +7.86%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '2.91'
$ws.Range('E33').Value = '  +3.55%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '2.27'
$ws.Range('E34').Value = '  +17.62%  '
$ws.Range('D35').Value = '151.25'
$ws.Range('E35').Value = '  +3.11%  '
$ws.Range('D36').Value = '0.0817'
$ws.Range('E36').Value = '  +5.55%  '
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  +2.01%  '
$ws.Range('D38').Value = '0.120'
$ws.Range('E38').Value = '  +4.02%  '
$ws.Range('D39').Value = '4.12'
$ws.Range('E39').Value = '  +5.22%  '
$ws.Range('D40').Value = '15.35'
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('D41').Value = '3.53'
$ws.Range('E41').Value = '  +9.01%  '
$ws.Range('D42').Value = '0.0317'
$ws.Range('E42').Value = '  +5.70%  '
$ws.Range('D43').Value = '2.034.88'
$ws.Range('E43').Value = '  +5.32%  '
$ws.Range('D44').Value = '18.85'
$ws.Range('E44').Value = '  +33.64%  '
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '89.71'
$ws.Range('E46').Value = '  -1.92%  '
$ws.Range('D47').Value = '9.09'
$ws.Range('E47').Value = '  +7.88%  '
$ws.Range('D48').Value = '108.15'
$ws.Range('E48').Value = '  +9.44%  '
$ws.Range('D49').Value = '1.74'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.839.55'
$ws.Range('E50').Value = '  +8.59%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.196'
$ws.Range('E51').Value = '  +5.06%  '
